$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Bug fix: the three "FINFISH,<5%FAT,SIMMERED,WO/DRIPPING" composite rows
# (retention-factor lookup, rows 60-62) were carrying three different/wrong
# retention-group codes (2402 / 2403 / 2404). They all belong to the same
# retention group as the row just above them (2401), so correct all three.
#
# Column B is formatted as Text (@), so a plain `.Value = 2401` would type
# the *string* "2401" into the cell (faithful Excel behaviour for entering a
# number into a Text-formatted cell). The source data stores these as real
# numbers, so briefly switch each cell to the Normal/General style while
# writing the numeric value, then restore the Text number format (and the
# top-aligned vertical alignment the column otherwise uses) so the cell
# keeps its original look.
foreach ($addr in @("B60", "B61", "B62")) {
    $cell = $ws.Range($addr)
    $cell.Style = "Normal"
    $cell.Value = 2401
    $cell.NumberFormat = "@"
    $cell.VerticalAlignment = -4160
}

# --- Update the sheet's saved selection to reflect where editing left off.
$ws.Activate()
$ws.Range("L58").Select()
